$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.28
$ws.Range("G2").Value = 2.44
$ws.Range("K2").Value = 3.9
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.28
$ws.Range("Q2").Value = 1.82
$ws.Range("R2").Value = 1.41
$ws.Range("S2").Value = 3.05
$ws.Range("T2").Value = 1.68
$ws.Range("U2").Value = 2.28
$ws.Range("W2").Value = 1.7
$ws.Range("X2").Value = 970
$ws.Range("Y2").Value = 970
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = 60
$ws.Range("AB2").Value = 11.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AD2").Value = 970
$ws.Range("AE2").Value = 36
$ws.Range("AF2").Value = 970
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 970
$ws.Range("AI2").Value = 44
$ws.Range("AJ2").Value = 32
$ws.Range("AK2").Value = 24
$ws.Range("AL2").Value = 36
$ws.Range("AM2").Value = 80
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2.2
$ws.Range("Q3").Value = 1.74
$ws.Range("S3").Value = 2.8
$ws.Range("T3").Value = 1.63
$ws.Range("U3").Value = 2.38
$ws.Range("X3").Value = 21
$ws.Range("AB3").Value = 15.5
$ws.Range("AJ3").Value = 1000
$ws.Range("G4").Value = 2.04
$ws.Range("H4").Value = 1.98
$ws.Range("H5").Value = 2.74
$ws.Range("I5").Value = 2.8
$ws.Range("K5").Value = 3.85
$ws.Range("R5").Value = 1.57
$ws.Range("AH5").Value = 15
$ws.Range("F6").Value = 3.95
$ws.Range("I6").Value = 2.04
$ws.Range("J6").Value = 3.85
$ws.Range("Q6").Value = 1.75
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 980
$ws.Range("F7").Value = 3.55
$ws.Range("G7").Value = 6.4
$ws.Range("I7").Value = 2.08
$ws.Range("J7").Value = 2.78
$ws.Range("F8").Value = 2.16
$ws.Range("H8").Value = 3.65
$ws.Range("I8").Value = 3.8
$ws.Range("J8").Value = 3.6
$ws.Range("K8").Value = 3.7
$ws.Range("P8").Value = 1.95
$ws.Range("R8").Value = 1.36
